$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a header row (row 1) and a single data row (row 2)
# describing the "SzVitalSigns" profile. We are adding a new profile
# ("SzCauseOfDeath") that needs three rows, and the existing SzVitalSigns
# data moves down to make room, ending up on row 5. Copy the existing row
# (values + formatting) down to row 5 first so nothing is lost, and copy
# its formatting into the two brand new rows (3 and 4) as well.
$ws.Range("A2:K2").Copy($ws.Range("A5:K5"))
$ws.Range("A2:K2").Copy($ws.Range("A3:K4"))

# Row 2: SzCauseOfDeath - first code (LOINC#79378-6)
$ws.Range("A2").Value = "SzCauseOfDeath"
$ws.Range("B2").Value = "Eswatini Cause of Death Profile"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "LOINC#79378-6"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "dateTime, Period, Timing, instant"
$ws.Range("H2").Value = "CodeableConcept"
$ws.Range("I2").Value = "optional"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""

# Row 3: SzCauseOfDeath - continuation (Value Types: string)
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "Eswatini Cause of Death Profile"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "string"
$ws.Range("I3").Value = "optional"
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""

# Row 4: SzCauseOfDeath - continuation (second code, LOINC#69440-6)
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "Eswatini Cause of Death Profile"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "LOINC#69440-6"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "Quantity, CodeableConcept, string"
$ws.Range("I4").Value = "optional"
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""

# Row 5: existing SzVitalSigns row (moved down) - update its Category Code
$ws.Range("C5").Value = "Observation Category Codes#vital-signs"
